$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.483.27'
$ws.Range('E2').Value = '  +1.84%  '
$ws.Range('D3').Value = '1.662.93'
$ws.Range('E3').Value = '  +1.13%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9988'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '235.58'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.89%  '
$ws.Range('E6').Value = '  +0.05%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4624'
$ws.Range('D7').Style = 'Normal'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2572'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -1.08%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06137'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.43%  '
$ws.Range('D10').Value = '1.661.50'
$ws.Range('E10').Value = '  +1.04%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.06943'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -1.51%  '
$ws.Range('E12').Value = '  -0.31%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.327'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.82%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '75.05'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.76%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.5730'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -3.75%  '
$ws.Range('E16').Value = '  +0.01%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.9996'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.02%  '
$ws.Range('E18').Value = '  +1.80%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000006703'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +1.57%  '
$ws.Range('E20').Value = '  +0.65%  '
$ws.Range('D21').Value = '1.872.44'
$ws.Range('E21').Value = '  +0.95%  '
$ws.Range('E22').Value = '  +1.69%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '8.643'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.63%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '5.211'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.54%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '134.71'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.06%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '14.87'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E27').Value = '  -1.69%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.707'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +3.54%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '103.58'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.25%  '
$ws.Range('E30').Value = '  +1.59%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.07703'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.18%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.591'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.60%  '
$ws.Range('E33').Value = '  +0.74%  '
$ws.Range('E34').Value = '  +1.69%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.5999'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.99%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.9387'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.99%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.9081'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +5.00%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.419'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -6.26%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '107.49'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +8.83%  '
$ws.Range('E40').Value = '  -0.03%  '
$ws.Range('B41').Value = 'VeChain'
$ws.Range('C41').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.01456'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -3.69%  '
$ws.Range('B42').Value = 'RenderToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.819'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +3.20%  '
$ws.Range('B43').Value = 'TheSandbox'
$ws.Range('C43').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.3706'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.06%  '
$ws.Range('B44').Value = 'FraxShare'
$ws.Range('C44').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '4.999'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +6.77%  '
$ws.Range('E45').Value = '  +0.44%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.05256'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.91%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '6.111'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.34%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '30.50'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +4.72%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.613'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +7.14%  '
$ws.Range('E50').Value = '  +0.17%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.9992'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.20%  '
